# Apply updated crypto market data (prices / 1h volume %) to Sheet1.
# Generated from the authoritative diff: only cells whose text actually
# changed are touched, each written as literal text (leading apostrophe)
# so Excel does not reinterpret numeric-looking strings like "28.409.32"
# or "0.9995" as numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.409.32'
$ws.Range("E2").Value = '''  -0.19%  '

$ws.Range("D3").Value = '''1.834.61'
$ws.Range("E3").Value = '''  +2.01%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '''  -0.07%  '

$ws.Range("D5").Value = '''318.47'
$ws.Range("E5").Value = '''  +0.59%  '

$ws.Range("D6").Value = '''0.9995'
$ws.Range("E6").Value = '''  -0.20%  '

$ws.Range("D7").Value = '''0.5312'
$ws.Range("E7").Value = '''  -2.03%  '

$ws.Range("D8").Value = '''0.3999'
$ws.Range("E8").Value = '''  +5.84%  '

$ws.Range("E9").Value = '''  +1.20%  '

$ws.Range("D10").Value = '''41.88'
$ws.Range("E10").Value = '''  -0.11%  '

$ws.Range("D11").Value = '''1.105'
$ws.Range("E11").Value = '''  -0.37%  '

$ws.Range("D12").Value = '''6.313'
$ws.Range("E12").Value = '''  +2.49%  '

$ws.Range("B13").Value = '''BinanceUSD'
$ws.Range("C13").Value = '''https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D13").Value = '''1.001'
$ws.Range("E13").Value = '''  -0.14%  '

$ws.Range("B14").Value = '''Chainlink'
$ws.Range("C14").Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''7.610'
$ws.Range("E14").Value = '''  +4.19%  '

$ws.Range("D15").Value = '''20.76'
$ws.Range("E15").Value = '''  +0.32%  '

$ws.Range("D16").Value = '''1.831.57'
$ws.Range("E16").Value = '''  +2.25%  '

$ws.Range("D17").Value = '''89.86'
$ws.Range("E17").Value = '''  +0.26%  '

$ws.Range("D18").Value = '''0.00001071'
$ws.Range("E18").Value = '''  +0.51%  '

$ws.Range("D19").Value = '''0.06596'
$ws.Range("E19").Value = '''  +1.25%  '

$ws.Range("D20").Value = '''17.61'
$ws.Range("E20").Value = '''  +0.98%  '

$ws.Range("D21").Value = '''0.9997'
$ws.Range("E21").Value = '''  -0.17%  '

$ws.Range("D22").Value = '''6.061'
$ws.Range("E22").Value = '''  +1.94%  '

$ws.Range("D23").Value = '''28.423.50'
$ws.Range("E23").Value = '''  -0.12%  '

$ws.Range("D24").Value = '''11.19'
$ws.Range("E24").Value = '''  +0.73%  '

$ws.Range("D25").Value = '''2.112'
$ws.Range("E25").Value = '''  +1.83%  '

$ws.Range("D26").Value = '''156.56'
$ws.Range("E26").Value = '''  -1.47%  '

$ws.Range("B27").Value = '''LidoDAOToken'
$ws.Range("C27").Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '''2.427'
$ws.Range("E27").Value = '''  +4.17%  '

$ws.Range("B28").Value = '''EthereumClassic'
$ws.Range("C28").Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''20.55'
$ws.Range("E28").Value = '''  +0.36%  '

$ws.Range("D29").Value = '''2.034.38'

$ws.Range("D30").Value = '''123.75'
$ws.Range("E30").Value = '''  +0.68%  '

$ws.Range("D31").Value = '''1.115'
$ws.Range("E31").Value = '''  +0.20%  '

$ws.Range("D32").Value = '''0.1098'
$ws.Range("E32").Value = '''  +3.88%  '

$ws.Range("E33").Value = '''  +1.27%  '

$ws.Range("D34").Value = '''5.638'
$ws.Range("E34").Value = '''  +0.34%  '

$ws.Range("D35").Value = '''0.07265'
$ws.Range("E35").Value = '''  +12.02%  '

$ws.Range("D36").Value = '''0.2249'
$ws.Range("E36").Value = '''  -1.10%  '

$ws.Range("B37").Value = '''VeChain'
$ws.Range("C37").Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.02341'
$ws.Range("E37").Value = '''  +1.95%  '

$ws.Range("B38").Value = '''InternetComputer(DFINITY)'
$ws.Range("C38").Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '''5.239'
$ws.Range("E38").Value = '''  +4.35%  '

$ws.Range("D39").Value = '''8.817'
$ws.Range("E39").Value = '''  +2.23%  '

$ws.Range("D40").Value = '''11.35'
$ws.Range("E40").Value = '''  +1.21%  '

$ws.Range("D41").Value = '''0.6264'
$ws.Range("E41").Value = '''  +0.89%  '

$ws.Range("D42").Value = '''1.198'
$ws.Range("E42").Value = '''  +0.38%  '

$ws.Range("D43").Value = '''1.414'
$ws.Range("E43").Value = '''  -2.62%  '

$ws.Range("B44").Value = '''Frax'
$ws.Range("C44").Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '''0.9988'
$ws.Range("E44").Value = '''  -0.22%  '

$ws.Range("B45").Value = '''EnergySwap'
$ws.Range("C45").Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''13.47'
$ws.Range("E45").Value = '''  +1.11%  '

$ws.Range("B46").Value = '''PancakeSwap'
$ws.Range("C46").Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '''3.704'
$ws.Range("E46").Value = '''  +0.47%  '

$ws.Range("B47").Value = '''Decentraland'
$ws.Range("C47").Value = '''https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.5823'
$ws.Range("E47").Value = '''  -0.17%  '

$ws.Range("B48").Value = '''Quant'
$ws.Range("C48").Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '''125.78'
$ws.Range("E48").Value = '''  -1.14%  '

$ws.Range("B49").Value = '''NEARProtocol'
$ws.Range("C49").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '''1.974'
$ws.Range("E49").Value = '''  +1.19%  '

$ws.Range("B50").Value = '''EOS'
$ws.Range("C50").Value = '''https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '''1.194'
$ws.Range("E50").Value = '''  -1.19%  '

$ws.Range("B51").Value = '''Cronos'
$ws.Range("C51").Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.06910'
$ws.Range("E51").Value = '''  +0.21%  '
